# Fix's and new Server
# - Updated DBs names to match the ones on the current server
# - Fixed bug #119 (duplicate/garbled "PRIORITY SERVICE" row whose description
#   actually belonged to a separate "Priority Service Hall Calls" feature)
# - Fixed bug #113 (missing special-feature catalog entries)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Remove the bad duplicate row: Name="PRIORITY SERVICE",
#        Description="Priority Service Hall Calls" ------------------------
$rowCount = $lo.ListRows.Count()
for ($i = $rowCount; $i -ge 1; $i--) {
    $lr = $lo.ListRows.Item($i)
    $rng = $lr.Range()
    $nameVal = $rng.Cells.Item(1, 1).Value()
    $descVal = $rng.Cells.Item(1, 2).Value()
    if ($nameVal -eq "PRIORITY SERVICE" -and $descVal -eq "Priority Service Hall Calls") {
        $rng.EntireRow.Delete()
    }
}

# --- 2. Append the new special-feature rows --------------------------------
$newFeatures = @(
    @("ADD-A-PHASE", "This job uses Ronk Add-A-Phase converter"),
    @("EMERGENCY COMMANDEERING", "This job has Emergency Commandeering features"),
    @("EMERGENCY DISPATCH", "This job monitors hall call bus for failures (wild operation)"),
    @("FIRE SERVICE ACCESS ELEVATOR", "This job is designated as a Fire Service Access Elevator (FSAE)"),
    @("GONG BOARD", "This job uses the Gong Board for Hall Lanterns"),
    @("OA BOARD", "This job uses the OA (output adapter) board"),
    @("PI: VIDATECH", "Vidatech Position Indicators")
)

foreach ($feature in $newFeatures) {
    $newRow = $lo.ListRows.Add()
    $rng = $newRow.Range()
    $rng.Cells.Item(1, 1).Value2 = $feature[0]
    $rng.Cells.Item(1, 2).Value2 = $feature[1]
}

# --- 3. Sort the whole table A-Z by the Name column -------------------------
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lo.ListColumns.Item(1).Range())
$lo.Sort.Header = 1
$lo.Sort.Apply()
